{"js": "// Replace each math-equation cell's text with its updated equation.\n// Mapping is old exact cell text -> new exact cell text, derived from the\n// authoritative OOXML diff. Every \"old\" string is unique in the document,\n// so an exact, case-sensitive whole-match search safely targets exactly one\n// cell each time.\nconst replacements = [\n  [\"20+78=98\", \"50+8=58\"],\n  [\"74-38=36\", \"98-67=31\"],\n  [\"99-51=48\", \"90+3=93\"],\n  [\"66-27=39\", \"10+6=16\"],\n  [\"53-27=26\", \"51-3=48\"],\n  [\"67-62=5\", \"11+49=60\"],\n  [\"98-24=74\", \"76-12=64\"],\n  [\"79-74=5\", \"69+7=76\"],\n  [\"79+19=98\", \"57-44=13\"],\n  [\"95-57=38\", \"71-29=42\"],\n  [\"13+82=95\", \"67+16=83\"],\n  [\"36-3=33\", \"70-13=57\"],\n  [\"68+22=90\", \"35+39=74\"],\n  [\"50-27=23\", \"25+13=38\"],\n  [\"79+10=89\", \"60-2=58\"],\n  [\"69+20=89\", \"31-0=31\"],\n  [\"78-11=67\", \"76-58=18\"],\n  [\"44-8=36\", \"71+23=94\"],\n  [\"91-12=79\", \"25+13=38\"],\n  [\"87-77=10\", \"32-14=18\"],\n  [\"5+57=62\", \"32+38=70\"],\n  [\"28+6=34\", \"56+32=88\"],\n  [\"49+15=64\", \"28-20=8\"],\n  [\"97-51=46\", \"55+10=65\"],\n  [\"22+50=72\", \"65-58=7\"],\n  [\"41+21=62\", \"68-45=23\"],\n  [\"35+5=40\", \"88-64=24\"],\n  [\"74-17=57\", \"48+48=96\"],\n  [\"13+20=33\", \"6+7=13\"],\n  [\"24+66=90\", \"75-44=31\"],\n  [\"16+8=24\", \"68-44=24\"],\n  [\"99-42=57\", \"50+11=61\"],\n  [\"84+2=86\", \"71-8=63\"],\n  [\"14+6=20\", \"11+50=61\"],\n  [\"30+64=94\", \"63-34=29\"],\n  [\"57+10=67\", \"82-26=56\"],\n  [\"18-1=17\", \"10+20=30\"],\n  [\"71-39=32\", \"74-7=67\"],\n  [\"48+4=52\", \"53-48=5\"],\n  [\"99-8=91\", \"93-55=38\"],\n  [\"81-7=74\", \"23-10=13\"],\n  [\"41+42=83\", \"77-62=15\"],\n  [\"24+63=87\", \"73-12=61\"],\n  [\"85-26=59\", \"24+61=85\"],\n  [\"13+5=18\", \"30-2=28\"],\n  [\"62-4=58\", \"35-27=8\"],\n  [\"17+48=65\", \"13+71=84\"],\n  [\"49-44=5\", \"80-42=38\"],\n  [\"92-31=61\", \"60+28=88\"],\n  [\"63-5=58\", \"62-49=13\"],\n  [\"91-19=72\", \"44-16=28\"],\n  [\"75-0=75\", \"87-53=34\"],\n  [\"37+18=55\", \"65-0=65\"],\n  [\"79-70=9\", \"78-60=18\"],\n  [\"44-29=15\", \"27+8=35\"],\n  [\"99-69=30\", \"1+26=27\"],\n  [\"65-1=64\", \"15-12=3\"],\n  [\"5+39=44\", \"16-14=2\"],\n  [\"6+31=37\", \"98-52=46\"],\n  [\"1+11=12\", \"6+7=13\"],\n  [\"50-1=49\", \"69-61=8\"],\n  [\"58-18=40\", \"46-32=14\"],\n  [\"47+26=73\", \"35-23=12\"],\n  [\"50+45=95\", \"44-1=43\"],\n  [\"28-27=1\", \"73+9=82\"],\n  [\"22+55=77\", \"18+59=77\"],\n  [\"6+4=10\", \"17+47=64\"],\n  [\"2+6=8\", \"39-10=29\"],\n  [\"12+47=59\", \"82-76=6\"],\n  [\"84-61=23\", \"37+44=81\"],\n  [\"72+14=86\", \"18+20=38\"],\n  [\"92-70=22\", \"87+2=89\"],\n  [\"28-22=6\", \"7+80=87\"],\n  [\"87-63=24\", \"45+31=76\"],\n  [\"30-17=13\", \"48+44=92\"],\n  [\"25+0=25\", \"78-55=23\"],\n  [\"44-4=40\", \"11+74=85\"],\n  [\"39+60=99\", \"1+91=92\"],\n  [\"92-37=55\", \"29+27=56\"],\n  [\"6+23=29\", \"38-4=34\"],\n  [\"37+58=95\", \"94-3=91\"],\n  [\"48-15=33\", \"81+18=99\"],\n  [\"13-12=1\", \"66-62=4\"],\n  [\"55-36=19\", \"12+14=26\"],\n  [\"0+50=50\", \"52-44=8\"],\n  [\"62+19=81\", \"93-67=26\"],\n  [\"4+13=17\", \"88-40=48\"],\n  [\"74-64=10\", \"79-13=66\"],\n  [\"37-25=12\", \"95-15=80\"],\n  [\"80-57=23\", \"73-62=11\"],\n  [\"15+14=29\", \"7+78=85\"],\n  [\"50-37=13\", \"7+41=48\"],\n  [\"3+64=67\", \"10+63=73\"],\n  [\"54-36=18\", \"40-10=30\"],\n  [\"61+12=73\", \"84-71=13\"],\n  [\"70+6=76\", \"47+37=84\"],\n  [\"21-9=12\", \"13+34=47\"],\n  [\"53-36=17\", \"95-24=71\"],\n  [\"48-11=37\", \"55-30=25\"],\n  [\"43+25=68\", \"85-37=48\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Each \"old\" equation text is unique in the document (verified against the\n  // source diff), so exactly one hit is expected; replace whichever range(s)\n  // search returns to stay robust either way.\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each math-equation cell's text with its updated equation using\n# Word's Find/Replace (Find.Execute), matching the authoritative OOXML diff.\n# Every \"old\" equation string is unique in the document, so an exact,\n# case-sensitive, whole-word match safely targets exactly one cell per pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"20+78=98\", \"50+8=58\"),\n    @(\"74-38=36\", \"98-67=31\"),\n    @(\"99-51=48\", \"90+3=93\"),\n    @(\"66-27=39\", \"10+6=16\"),\n    @(\"53-27=26\", \"51-3=48\"),\n    @(\"67-62=5\", \"11+49=60\"),\n    @(\"98-24=74\", \"76-12=64\"),\n    @(\"79-74=5\", \"69+7=76\"),\n    @(\"79+19=98\", \"57-44=13\"),\n    @(\"95-57=38\", \"71-29=42\"),\n    @(\"13+82=95\", \"67+16=83\"),\n    @(\"36-3=33\", \"70-13=57\"),\n    @(\"68+22=90\", \"35+39=74\"),\n    @(\"50-27=23\", \"25+13=38\"),\n    @(\"79+10=89\", \"60-2=58\"),\n    @(\"69+20=89\", \"31-0=31\"),\n    @(\"78-11=67\", \"76-58=18\"),\n    @(\"44-8=36\", \"71+23=94\"),\n    @(\"91-12=79\", \"25+13=38\"),\n    @(\"87-77=10\", \"32-14=18\"),\n    @(\"5+57=62\", \"32+38=70\"),\n    @(\"28+6=34\", \"56+32=88\"),\n    @(\"49+15=64\", \"28-20=8\"),\n    @(\"97-51=46\", \"55+10=65\"),\n    @(\"22+50=72\", \"65-58=7\"),\n    @(\"41+21=62\", \"68-45=23\"),\n    @(\"35+5=40\", \"88-64=24\"),\n    @(\"74-17=57\", \"48+48=96\"),\n    @(\"13+20=33\", \"6+7=13\"),\n    @(\"24+66=90\", \"75-44=31\"),\n    @(\"16+8=24\", \"68-44=24\"),\n    @(\"99-42=57\", \"50+11=61\"),\n    @(\"84+2=86\", \"71-8=63\"),\n    @(\"14+6=20\", \"11+50=61\"),\n    @(\"30+64=94\", \"63-34=29\"),\n    @(\"57+10=67\", \"82-26=56\"),\n    @(\"18-1=17\", \"10+20=30\"),\n    @(\"71-39=32\", \"74-7=67\"),\n    @(\"48+4=52\", \"53-48=5\"),\n    @(\"99-8=91\", \"93-55=38\"),\n    @(\"81-7=74\", \"23-10=13\"),\n    @(\"41+42=83\", \"77-62=15\"),\n    @(\"24+63=87\", \"73-12=61\"),\n    @(\"85-26=59\", \"24+61=85\"),\n    @(\"13+5=18\", \"30-2=28\"),\n    @(\"62-4=58\", \"35-27=8\"),\n    @(\"17+48=65\", \"13+71=84\"),\n    @(\"49-44=5\", \"80-42=38\"),\n    @(\"92-31=61\", \"60+28=88\"),\n    @(\"63-5=58\", \"62-49=13\"),\n    @(\"91-19=72\", \"44-16=28\"),\n    @(\"75-0=75\", \"87-53=34\"),\n    @(\"37+18=55\", \"65-0=65\"),\n    @(\"79-70=9\", \"78-60=18\"),\n    @(\"44-29=15\", \"27+8=35\"),\n    @(\"99-69=30\", \"1+26=27\"),\n    @(\"65-1=64\", \"15-12=3\"),\n    @(\"5+39=44\", \"16-14=2\"),\n    @(\"6+31=37\", \"98-52=46\"),\n    @(\"1+11=12\", \"6+7=13\"),\n    @(\"50-1=49\", \"69-61=8\"),\n    @(\"58-18=40\", \"46-32=14\"),\n    @(\"47+26=73\", \"35-23=12\"),\n    @(\"50+45=95\", \"44-1=43\"),\n    @(\"28-27=1\", \"73+9=82\"),\n    @(\"22+55=77\", \"18+59=77\"),\n    @(\"6+4=10\", \"17+47=64\"),\n    @(\"2+6=8\", \"39-10=29\"),\n    @(\"12+47=59\", \"82-76=6\"),\n    @(\"84-61=23\", \"37+44=81\"),\n    @(\"72+14=86\", \"18+20=38\"),\n    @(\"92-70=22\", \"87+2=89\"),\n    @(\"28-22=6\", \"7+80=87\"),\n    @(\"87-63=24\", \"45+31=76\"),\n    @(\"30-17=13\", \"48+44=92\"),\n    @(\"25+0=25\", \"78-55=23\"),\n    @(\"44-4=40\", \"11+74=85\"),\n    @(\"39+60=99\", \"1+91=92\"),\n    @(\"92-37=55\", \"29+27=56\"),\n    @(\"6+23=29\", \"38-4=34\"),\n    @(\"37+58=95\", \"94-3=91\"),\n    @(\"48-15=33\", \"81+18=99\"),\n    @(\"13-12=1\", \"66-62=4\"),\n    @(\"55-36=19\", \"12+14=26\"),\n    @(\"0+50=50\", \"52-44=8\"),\n    @(\"62+19=81\", \"93-67=26\"),\n    @(\"4+13=17\", \"88-40=48\"),\n    @(\"74-64=10\", \"79-13=66\"),\n    @(\"37-25=12\", \"95-15=80\"),\n    @(\"80-57=23\", \"73-62=11\"),\n    @(\"15+14=29\", \"7+78=85\"),\n    @(\"50-37=13\", \"7+41=48\"),\n    @(\"3+64=67\", \"10+63=73\"),\n    @(\"54-36=18\", \"40-10=30\"),\n    @(\"61+12=73\", \"84-71=13\"),\n    @(\"70+6=76\", \"47+37=84\"),\n    @(\"21-9=12\", \"13+34=47\"),\n    @(\"53-36=17\", \"95-24=71\"),\n    @(\"48-11=37\", \"55-30=25\"),\n    @(\"43+25=68\", \"85-37=48\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
